# Commit: "added new fields to output"
# Added date_different between date of collection and date of travel, and
# duration as well. In this metadata sheet this shows up as a third route
# of data (Route3 / EUS,LIV <-> LIV,EUS) being appended in column D, the
# existing "Route" column (now "Route2") being re-pointed at the CMD/SRA
# route's trimmed sample times, and the original "Route 1" column being
# trimmed down to shorter representative time samples.
#
# Cells that already held a value before this edit use a leading apostrophe
# so the engine treats the new value as literal text (not a number) and so
# the pre-existing cell style (quote-prefix / number-format) carries over
# unchanged, matching how Excel itself preserves a cell's format when you
# retype its contents. Brand new cells (D1:D3) are written without the
# apostrophe since they start out with the workbook's default (no) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header ---
$ws.Range("C1").Value = "Route2"
$ws.Range("D1").Value = "Route3"

# --- Row 2: origin and destination down ---
$ws.Range("C2").Value = "'CMD,SRA"
$ws.Range("D2").Value = "EUS,LIV"

# --- Row 3: origin and destination up ---
$ws.Range("C3").Value = "'SRA,CMD"
$ws.Range("D3").Value = "LIV,EUS"

# --- Row 4: downweekdaytime ---
$ws.Range("B4").Value = "'0612,0700"
$ws.Range("C4").Value = "'0800,0900"
$ws.Range("D4").Value = "'1500,1700,1800"

# --- Row 5: downsaturdaytimes ---
$ws.Range("B5").Value = "'0612"
$ws.Range("C5").Value = "'1100"
$ws.Range("D5").Value = "'0900,1000,1100,1200"

# --- Row 6: downsundaytimes ---
$ws.Range("B6").Value = "'0848"
$ws.Range("C6").Value = "'1600"
$ws.Range("D6").Value = "'1300,1500,1612"

# --- Row 7: upweekdaytimes ---
$ws.Range("B7").Value = "'0656"
$ws.Range("C7").Value = "'1700"
$ws.Range("D7").Value = "'1500,1700,1800"
# D7 moves onto the same quote-prefixed/number-format style already used by D4.
$ws.Range("D7").NumberFormat = $ws.Range("D4").NumberFormat

# --- Row 8: upsaturdaytimes ---
$ws.Range("B8").Value = "'0626"
$ws.Range("C8").Value = "'2000"
$ws.Range("D8").Value = "'0900,1000,1100,1200"

# --- Row 9: upsundaytimes ---
$ws.Range("B9").Value = "'1030"
$ws.Range("C9").Value = "'1900"
$ws.Range("D9").Value = "'1300,1500,1612"

# Move the active selection to D12, matching where the editor ended up after
# entering the new data.
$ws.Range("D12").Select()
